$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 39 ("Boost broken?"): content placeholder title line gets
# re-cased/re-worded ("A remedy for the review manager starvation." ->
# "A Remedy for the Review Manager Starvation.") split across several
# runs, and the paragraph is centered with bullets turned off.
# ---------------------------------------------------------------------
$s39   = $p.Slides.Item(39)
$shp39 = $s39.Shapes.Item(5)
$tr39  = $shp39.TextFrame.TextRange

$tr39.Text = "A Remedy for "
$tr39.ParagraphFormat.Alignment = 2      # ppAlignCenter
$tr39.ParagraphFormat.Bullet.Visible = 0 # a:buNone

$run = $tr39.InsertAfter("the Review ")
$run = $run.InsertAfter("M")
$run = $run.InsertAfter("anager ")
$run = $run.InsertAfter("S")
$run = $run.InsertAfter("tarvation")
$run = $run.InsertAfter(".")

# ---------------------------------------------------------------------
# Slide 40 ("Where is the motivation?"): the word "approval" in the
# 5th bullet paragraph is highlighted in blue (0070C0).
# ---------------------------------------------------------------------
$s40   = $p.Slides.Item(40)
$shp40 = $s40.Shapes.Item(5)
$tr40  = $shp40.TextFrame.TextRange

$para  = $tr40.Paragraphs(5)
$word  = $para.Characters(21, 8)   # "approval"
$word.Font.Color.RGB = 12611584    # RGB(0, 112, 192) = 0070C0

# ---------------------------------------------------------------------
# Slide 42 ("A Review Manager Assistant"): the two runs making up the
# tail of the second paragraph ("oes all the work ... file a " and
# "final report.") are merged back into a single run.
# ---------------------------------------------------------------------
$s42   = $p.Slides.Item(42)
$shp42 = $s42.Shapes.Item(5)
$tr42  = $shp42.TextFrame.TextRange

$para2 = $tr42.Paragraphs(2)
$runA  = $para2.Runs(2)
$runB  = $para2.Runs(3)
$runA.Text = $runA.Text + $runB.Text
$runB.Text = ""
